# Insert a new data row at row 243 (pushing existing rows 243:300 down to 244:301)
# and populate it with a new weekly price record for Acelga.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(243).Insert()

$ws.Cells.Item(243, 1).Value = 5
$ws.Cells.Item(243, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(243, 3).Value = "Maule"
$ws.Cells.Item(243, 4).Value = 44798
$ws.Cells.Item(243, 5).Value = 7
$ws.Cells.Item(243, 6).Value = 100112009
$ws.Cells.Item(243, 7).Value = "Acelga"
$ws.Cells.Item(243, 8).Value = "Sin especificar"
$ws.Cells.Item(243, 9).Value = "Primera"
$ws.Cells.Item(243, 10).Value = 500
$ws.Cells.Item(243, 11).Value = 2500
$ws.Cells.Item(243, 12).Value = 2500
$ws.Cells.Item(243, 13).Value = 2500
$ws.Cells.Item(243, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(243, 15).Value = "Región del Maule"
$ws.Cells.Item(243, 16).Value = 625
$ws.Cells.Item(243, 17).Value = 4
$ws.Cells.Item(243, 18).Value = "Hortaliza"
